$d = $word.ActiveDocument

# Locate the first of the trailing empty paragraphs (the one right after the
# "Deny access ... regulations." paragraph) by searching for the last
# non-empty paragraph's following blank paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -match "not conform to regulations") {
        $target = $i + 1
        break
    }
}

$anchorPara = $d.Paragraphs.Item($target)
$anchorRange = $anchorPara.Range

# Make room: insert two new blank paragraphs right after the anchor blank
# paragraph, so the anchor paragraph itself is left untouched.
$null = $anchorRange.InsertParagraphAfter()
$firstNew = $d.Paragraphs.Item($target + 1)
$null = $firstNew.Range.InsertParagraphAfter()
$secondNew = $d.Paragraphs.Item($target + 2)

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml1 = "<w:p $ns><w:r><w:t>App development plug-in that encrypts user selected variables for them.</w:t></w:r></w:p>"
$null = $firstNew.Range.InsertXML($xml1)

$xml2 = "<w:p $ns>" + `
    "<w:r><w:t xml:space='preserve'>a nice front-end, scan the code for variables, and give the option to select specific variables from that list, and encrypt  with a selected method of encryption, with an optional  explanation. That way the user doesn't need to understand how to implement it themselves, they can encrypt data with three clicks. And, if we make it </w:t></w:r>" + `
    "<w:proofErr w:type='gramStart'/>" + `
    "<w:r><w:t>open-source</w:t></w:r>" + `
    "<w:proofErr w:type='gramEnd'/>" + `
    "<w:r><w:t xml:space='preserve'> the code can be improved upon by others / make the code more trustworthy.</w:t></w:r>" + `
    "</w:p>"
$null = $secondNew.Range.InsertXML($xml2)

Write-Output "Inserted two new idea paragraphs."
